$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 217, shifting existing rows 217..321 down to 218..322
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly data point
$ws.Range("A217").Value = 6
$ws.Range("B217").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C217").Value = "Metropolitana"
$ws.Range("D217").Value = 45089
$ws.Range("E217").Value = 13
$ws.Range("F217").Value = 100112001
$ws.Range("G217").Value = "Berenjena"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 400
$ws.Range("K217").Value = 5000
$ws.Range("L217").Value = 6000
$ws.Range("M217").Value = 5425
$ws.Range("N217").Value = "$/caja 50 unidades"
$ws.Range("O217").Value = "Región de Arica y Parinacota"
$ws.Range("P217").Value = 108
$ws.Range("Q217").Value = 50
$ws.Range("R217").Value = "Hortaliza"
